# Update the "Chart" sheet's rolling date-window export:
#  - drop the oldest date row (2025-10-24) and shift every other row up by one
#  - append a new row for the next day (2026-01-21) with its own "Items" count
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest data row (row 2); rows below shift up automatically,
# which also drops the now-unused "2025-10-24" shared string.
$ws.Rows.Item(2).Delete()

# The data block now occupies rows 2-89; populate the freed row 90 with the
# newest date. Writing the date through a formula and then pasting the
# result back as a value keeps it a plain text/shared-string cell instead of
# letting Excel auto-convert the literal into a date serial number.
$ws.Range("A90").Formula = '="2026-01-21"'
$ws.Range("A90").Copy()
$ws.Range("A90").PasteSpecial(-4163)

$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 24
